$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: these cells held their numbers as text; convert them to real numbers.
$ws.Range("A1").Value = 10.26
$ws.Range("A2").Value = 20.12
$ws.Range("A3").Value = 30.01
$ws.Range("A4").Value = 40.29
$ws.Range("A5").Value = 50.18

# Column B: used to hold empty-string placeholders in every row; clear them out,
# keeping the (now blank) cells present, and record the uploaded file name in B1.
$ws.Range("B1:B5").ClearContents()
$ws.Range("B2:B5").NumberFormat = "General"
$ws.Range("B1").Value = "New File"
